$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append an emoji to each team name in column A (rows 2-9), keeping the
# rest of the row (dollar figures / notes) untouched.
$ws.Range("A2").Value = "Otistics 🤥"
$ws.Range("A3").Value = "Los Yahoo 🍯"
$ws.Range("A4").Value = "Out of PO 🕋"
$ws.Range("A5").Value = "Team of Outs 🧙🏼"
$ws.Range("A6").Value = "MaltaSpor 🫄🏿"
$ws.Range("A7").Value = "NSY 🇸🇾"
$ws.Range("A8").Value = "TrendyOwls 🦉"
$ws.Range("A9").Value = "Orthopedics United 🏥"

# Restore the selection to A2 (matches the saved view state in the diff).
$ws.Range("A2").Select()
